$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 19) to the products table
$row = 19

$ws.Cells.Item($row, 1).Value = "yes sir"
$ws.Cells.Item($row, 2).Value = "test"
$ws.Cells.Item($row, 3).Value = "Don Valley"
$ws.Cells.Item($row, 4).Value = 12
$ws.Cells.Item($row, 5).Value = 15

# Expiry date - copy the date style used by the other rows (e.g. F18) and set the value
$ws.Cells.Item($row, 6).Value = (Get-Date -Year 2024 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item($row, 6).NumberFormat = "yyyy-mm-dd"

$ws.Cells.Item($row, 7).Value = "Anxiolytic"
$ws.Cells.Item($row, 8).Value = "Tablet"
